# Updated cryptos list - applies per-row Coin/Link/Price/Volume(1h) updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '31.292.31'
$ws.Range("E2").Value = '  +1.46%  '

# Row 3
$ws.Range("D3").Value = '2.000.38'
$ws.Range("E3").Value = '  +4.03%  '

# Row 4
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = '  +0.98%  '

# Row 5
$ws.Range("D5").Value = "'0.7751"
$ws.Range("E5").Value = '  +28.33%  '

# Row 6
$ws.Range("D6").Value = "'256.07"
$ws.Range("E6").Value = '  +2.09%  '

# Row 7
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = '  +0.74%  '

# Row 8
$ws.Range("D8").Value = "'0.3475"
$ws.Range("E8").Value = '  +14.31%  '

# Row 9
$ws.Range("D9").Value = "'28.10"
$ws.Range("E9").Value = '  +16.95%  '

# Row 10
$ws.Range("D10").Value = "'0.07178"
$ws.Range("E10").Value = '  +7.78%  '

# Row 11
$ws.Range("D11").Value = "'0.8467"
$ws.Range("E11").Value = '  +6.71%  '

# Row 12
$ws.Range("D12").Value = "'0.08198"
$ws.Range("E12").Value = '  +3.75%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '2.003.29'
$ws.Range("E13").Value = '  +4.47%  '

# Row 14
$ws.Range("B14").Value = 'Litecoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D14").Value = "'101.00"
$ws.Range("E14").Value = '  -0.71%  '

# Row 15
$ws.Range("D15").Value = "'5.651"
$ws.Range("E15").Value = '  +6.01%  '

# Row 16
$ws.Range("D16").Value = "'15.61"
$ws.Range("E16").Value = '  +15.64%  '

# Row 17
$ws.Range("D17").Value = "'272.50"
$ws.Range("E17").Value = '  -4.60%  '

# Row 18
$ws.Range("D18").Value = '31.291.17'
$ws.Range("E18").Value = '  +2.04%  '

# Row 19
$ws.Range("D19").Value = "'0.000008326"
$ws.Range("E19").Value = '  +9.13%  '

# Row 20
$ws.Range("D20").Value = "'6.006"
$ws.Range("E20").Value = '  +10.09%  '

# Row 21
$ws.Range("D21").Value = '2.263.85'

# Row 22
$ws.Range("D22").Value = "'0.9999"
$ws.Range("E22").Value = '  +0.36%  '

# Row 23
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = '  +1.06%  '

# Row 24
$ws.Range("D24").Value = "'7.137"
$ws.Range("E24").Value = '  +8.38%  '

# Row 25
$ws.Range("E25").Value = '  +8.66%  '

# Row 26
$ws.Range("D26").Value = "'164.77"
$ws.Range("E26").Value = '  +0.36%  '

# Row 27
$ws.Range("D27").Value = "'0.1416"
$ws.Range("E27").Value = '  +30.23%  '

# Row 28
$ws.Range("D28").Value = "'20.01"
$ws.Range("E28").Value = '  +2.58%  '

# Row 29
$ws.Range("D29").Value = "'2.418"
$ws.Range("E29").Value = '  +23.66%  '

# Row 30
$ws.Range("D30").Value = "'1.607"
$ws.Range("E30").Value = '  +5.43%  '

# Row 31
$ws.Range("D31").Value = "'4.653"
$ws.Range("E31").Value = '  +6.00%  '

# Row 32
$ws.Range("E32").Value = '  +1.16%  '

# Row 33
$ws.Range("D33").Value = "'4.487"
$ws.Range("E33").Value = '  +5.01%  '

# Row 34
$ws.Range("D34").Value = "'0.05376"
$ws.Range("E34").Value = '  +8.81%  '

# Row 35
$ws.Range("D35").Value = "'1.267"
$ws.Range("E35").Value = '  +9.22%  '

# Row 36
$ws.Range("D36").Value = "'0.7881"
$ws.Range("E36").Value = '  +10.18%  '

# Row 37
$ws.Range("D37").Value = "'2.778"
$ws.Range("E37").Value = '  -0.10%  '

# Row 38
$ws.Range("D38").Value = "'0.9992"
$ws.Range("E38").Value = '  +0.39%  '

# Row 39
$ws.Range("D39").Value = "'0.02014"
$ws.Range("E39").Value = '  +4.43%  '

# Row 40
$ws.Range("D40").Value = "'2.941"
$ws.Range("E40").Value = '  +1.26%  '

# Row 41
$ws.Range("D41").Value = "'86.20"
$ws.Range("E41").Value = '  +11.75%  '

# Row 42
$ws.Range("D42").Value = "'6.838"
$ws.Range("E42").Value = '  +7.08%  '

# Row 43
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").Value = "'0.4702"
$ws.Range("E43").Value = '  +6.27%  '

# Row 44
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = "'2.143"
$ws.Range("E44").Value = '  +6.64%  '

# Row 45
$ws.Range("D45").Value = "'0.8604"
$ws.Range("E45").Value = '  +1.65%  '

# Row 46
$ws.Range("D46").Value = "'105.52"
$ws.Range("E46").Value = '  +3.91%  '

# Row 47
$ws.Range("D47").Value = "'10.20"
$ws.Range("E47").Value = '  +1.05%  '

# Row 48
$ws.Range("D48").Value = "'0.9998"
$ws.Range("E48").Value = '  +0.37%  '

# Row 49
$ws.Range("D49").Value = "'7.784"
$ws.Range("E49").Value = '  +8.03%  '

# Row 50
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").Value = "'37.92"
$ws.Range("E50").Value = '  +6.60%  '

# Row 51
$ws.Range("B51").Value = 'SynthetixNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D51").Value = "'3.028"
$ws.Range("E51").Value = '  +42.81%  '

